$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values for rows 2 and 3
$ws.Range("B2").Value = 87
$ws.Range("B3").Value = 52

# Row 4 now holds what used to be row 5's A value (2) with an updated B value (25)
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 25

# Remove the old row 5 entirely (shifts rows up, shrinking the used range)
$ws.Rows("5:5").Delete()
